$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 797.5
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 797.5
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 797.5
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -1449.5
$ws.Range("H92").Value = 143
$ws.Range("J92").Value = 212.5
$ws.Range("L92").Value = 212.5
$ws.Range("N92").Value = -2708.5
$ws.Range("H94").Value = 8826.182000000001
$ws.Range("I94").Value = 8826.182000000001
$ws.Range("K94").Value = 8826.182000000001
$ws.Range("M94").Value = -8375.182000000001
$ws.Range("H100").Value = 7560.75
$ws.Range("I100").Value = 7243
$ws.Range("J100").Value = 7666.6665
$ws.Range("K100").Value = 7243
$ws.Range("L100").Value = 7666.6665
$ws.Range("M100").Value = -6702
$ws.Range("N100").Value = -8748.666499999999
$ws.Range("H138").Value = 2921.4707
$ws.Range("J138").Value = 3205.077
$ws.Range("L138").Value = 9615.231
$ws.Range("N138").Value = -19895.231
$ws.Range("H141").Value = 1842.1428
$ws.Range("I141").Value = 1044.5454
$ws.Range("K141").Value = 3133.6362
$ws.Range("M141").Value = 2046.3638
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13587.782
$ws.Range("I32").Value = 11150.277
$ws.Range("K32").Value = 11150.277
$ws.Range("M32").Value = -10863.277
$ws.Range("H61").Value = 4316
$ws.Range("I61").Value = 4479.2
$ws.Range("K61").Value = 4479.2
$ws.Range("M61").Value = -4267.2
$ws.Range("H74").Value = 2319.8
$ws.Range("I74").Value = 1868.5625
$ws.Range("J74").Value = 4124.75
$ws.Range("K74").Value = 1868.5625
$ws.Range("L74").Value = 4124.75
$ws.Range("M74").Value = -994.5625
$ws.Range("N74").Value = -5872.75
$ws.Range("H77").Value = 2319.8
$ws.Range("I77").Value = 1868.5625
$ws.Range("J77").Value = 4124.75
$ws.Range("K77").Value = 9342.8125
$ws.Range("L77").Value = 20623.75
$ws.Range("M77").Value = -4974.8125
$ws.Range("N77").Value = -29359.75
$ws.Range("H88").Value = 1549.3334
$ws.Range("I88").Value = 1598.3334
$ws.Range("J88").Value = 1516.6666
$ws.Range("K88").Value = 1598.3334
$ws.Range("L88").Value = 1516.6666
$ws.Range("M88").Value = -1192.3334
$ws.Range("N88").Value = -2328.6666
$ws.Range("H91").Value = 1549.3334
$ws.Range("I91").Value = 1598.3334
$ws.Range("J91").Value = 1516.6666
$ws.Range("K91").Value = 1598.3334
$ws.Range("L91").Value = 1516.6666
$ws.Range("M91").Value = -194.3334
$ws.Range("N91").Value = -4324.6666
$ws.Range("H132").Value = 7654
$ws.Range("I132").Value = 7654
$ws.Range("K132").Value = 22962
$ws.Range("M132").Value = -20432
$ws.Range("H136").Value = 4316
$ws.Range("I136").Value = 4479.2
$ws.Range("K136").Value = 13437.6
$ws.Range("M136").Value = -10887.6
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4470.077
$ws.Range("J86").Value = 6544.222
$ws.Range("L86").Value = 6544.222
$ws.Range("N86").Value = -8790.222
$ws.Range("H89").Value = 4470.077
$ws.Range("J89").Value = 6544.222
$ws.Range("L89").Value = 32721.11
$ws.Range("N89").Value = -43953.11
$ws.Range("H94").Value = 1100
$ws.Range("I94").Value = 1000
$ws.Range("K94").Value = 1000
$ws.Range("M94").Value = -549
$ws.Range("H99").Value = 3454.7334
$ws.Range("I99").Value = 3537.4546
$ws.Range("J99").Value = 3227.25
$ws.Range("K99").Value = 3537.4546
$ws.Range("L99").Value = 3227.25
$ws.Range("M99").Value = -2039.4546
$ws.Range("N99").Value = -6223.25
$ws.Range("H134").Value = 2370.4285
$ws.Range("I134").Value = 1099
$ws.Range("J134").Value = 9999
$ws.Range("K134").Value = 3297
$ws.Range("L134").Value = 29997
$ws.Range("M134").Value = -762
$ws.Range("N134").Value = -35067
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5395.7
$ws.Range("I31").Value = 4573.636
$ws.Range("J31").Value = 7656.375
$ws.Range("K31").Value = 4573.636
$ws.Range("L31").Value = 7656.375
$ws.Range("M31").Value = -4278.636
$ws.Range("N31").Value = -8246.375
$ws.Range("H34").Value = 5395.7
$ws.Range("I34").Value = 4573.636
$ws.Range("J34").Value = 7656.375
$ws.Range("K34").Value = 4573.636
$ws.Range("L34").Value = 7656.375
$ws.Range("M34").Value = -4371.636
$ws.Range("N34").Value = -8060.375
$ws.Range("H52").Value = 293748.75
$ws.Range("J52").Value = 375000
$ws.Range("L52").Value = 375000
$ws.Range("N52").Value = -375588
$ws.Range("H58").Value = 2435.2666
$ws.Range("I58").Value = 1925.4615
$ws.Range("J58").Value = 5749
$ws.Range("K58").Value = 1925.4615
$ws.Range("L58").Value = 5749
$ws.Range("M58").Value = -1722.4615
$ws.Range("N58").Value = -6155
$ws.Range("H132").Value = 2754.95
$ws.Range("I132").Value = 2000.5294
$ws.Range("K132").Value = 6001.5882
$ws.Range("M132").Value = -3471.5882
$ws.Range("H134").Value = 1965.0555
$ws.Range("I134").Value = 1799
$ws.Range("J134").Value = 7777
$ws.Range("K134").Value = 5397
$ws.Range("L134").Value = 23331
$ws.Range("M134").Value = -2862
$ws.Range("N134").Value = -28401
$ws.Range("H136").Value = 2435.2666
$ws.Range("I136").Value = 1925.4615
$ws.Range("J136").Value = 5749
$ws.Range("K136").Value = 5776.3845
$ws.Range("L136").Value = 17247
$ws.Range("M136").Value = -3226.3845
$ws.Range("N136").Value = -22347
$ws.Range("H138").Value = 74250
$ws.Range("J138").Value = 83333.336
$ws.Range("L138").Value = 83333.336
$ws.Range("N138").Value = -93613.336
$ws.Range("H139").Value = 50000
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 5997.5
$ws.Range("I3").Value = 5997.5
$ws.Range("K3").Value = 17992.5
$ws.Range("M3").Value = -17880.5
$ws.Range("H18").Value = 2750
$ws.Range("I18").Value = 500
$ws.Range("J18").Value = 5000
$ws.Range("K18").Value = 1500
$ws.Range("L18").Value = 15000
$ws.Range("M18").Value = -1331
$ws.Range("N18").Value = -15338
$ws.Range("H50").Value = 1300
$ws.Range("I50").Value = 100
$ws.Range("K50").Value = 300
$ws.Range("M50").Value = 181
$ws.Range("H53").Value = 1300
$ws.Range("I53").Value = 100
$ws.Range("K53").Value = 300
$ws.Range("M53").Value = 181
$ws.Range("H113").Value = 2142.6428
$ws.Range("J113").Value = 2562.8
$ws.Range("L113").Value = 7688.400000000001
$ws.Range("N113").Value = -12028.4
$ws.Range("H129").Value = 1816.6666
$ws.Range("I129").Value = 1133.3334
$ws.Range("J129").Value = 2500
$ws.Range("K129").Value = 3400.0002
$ws.Range("L129").Value = 7500
$ws.Range("M129").Value = 1599.9998
$ws.Range("N129").Value = -17500
$ws.Range("H132").Value = 3337.5417
$ws.Range("I132").Value = 3750.125
$ws.Range("J132").Value = 3131.25
$ws.Range("K132").Value = 33751.125
$ws.Range("L132").Value = 28181.25
$ws.Range("M132").Value = -31221.125
$ws.Range("N132").Value = -33241.25
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2247
$ws.Range("J80").Value = 2000
$ws.Range("L80").Value = 2000
$ws.Range("N80").Value = -3996
$ws.Range("H83").Value = 2247
$ws.Range("J83").Value = 2000
$ws.Range("L83").Value = 10000
$ws.Range("N83").Value = -19984
$ws.Range("H132").Value = 62643.527
$ws.Range("I132").Value = 72514.31
$ws.Range("J132").Value = 9999.333000000001
$ws.Range("K132").Value = 217542.93
$ws.Range("L132").Value = 29997.999
$ws.Range("M132").Value = -215012.93
$ws.Range("N132").Value = -35057.999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5933.846
$ws.Range("I7").Value = 4887.875
$ws.Range("K7").Value = 4887.875
$ws.Range("M7").Value = -4775.875
$ws.Range("H22").Value = 998.4666999999999
$ws.Range("J22").Value = 1193.75
$ws.Range("L22").Value = 1193.75
$ws.Range("N22").Value = -1783.75
$ws.Range("H27").Value = 998.4666999999999
$ws.Range("J27").Value = 1193.75
$ws.Range("L27").Value = 1193.75
$ws.Range("N27").Value = -1407.75
$ws.Range("H93").Value = 1825.7693
$ws.Range("I93").Value = 1849.5
$ws.Range("J93").Value = 1815.2222
$ws.Range("K93").Value = 1849.5
$ws.Range("L93").Value = 1815.2222
$ws.Range("M93").Value = -601.5
$ws.Range("N93").Value = -4311.2222
$ws.Range("H126").Value = 5933.846
$ws.Range("I126").Value = 4887.875
$ws.Range("K126").Value = 14663.625
$ws.Range("M126").Value = -12193.625
$ws.Range("H132").Value = 6399
$ws.Range("I132").Value = 2000
$ws.Range("K132").Value = 6000
$ws.Range("M132").Value = -3470
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1335.375
$ws.Range("J96").Value = 1349
$ws.Range("L96").Value = 1349
$ws.Range("N96").Value = -4095
$ws.Range("H100").Value = 608.5
$ws.Range("I100").Value = 429.41666
$ws.Range("K100").Value = 858.83332
$ws.Range("M100").Value = -317.83332
$ws.Range("H122").Value = 2993.3428
$ws.Range("I122").Value = 2113.318
$ws.Range("K122").Value = 6339.954000000001
$ws.Range("M122").Value = -3889.954000000001
$ws.Range("H132").Value = 4463.2354
$ws.Range("I132").Value = 4289.52
$ws.Range("J132").Value = 4945.778
$ws.Range("K132").Value = 12868.56
$ws.Range("L132").Value = 14837.334
$ws.Range("M132").Value = -10338.56
$ws.Range("N132").Value = -19897.334
$ws.Range("H136").Value = 5486.893
$ws.Range("I136").Value = 3875.1177
$ws.Range("J136").Value = 7977.8184
$ws.Range("K136").Value = 11625.3531
$ws.Range("L136").Value = 23933.4552
$ws.Range("M136").Value = -9075.3531
$ws.Range("N136").Value = -29033.4552
